# Update 2016 logins/passwords to 2017 with new generated passwords.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "stashevskiy.a.s.2017"
$ws.Range("E1").Value = "vnlCts78"

$ws.Range("D2").Value = "trubicin.yu.a.2017"
$ws.Range("E2").Value = "9wvEF7OT"

$ws.Range("D3").Value = "smirnov.s.v.2017"
$ws.Range("E3").Value = "8DMIhIyK"

$ws.Range("D4").Value = "gorohov.n.s.2017"
$ws.Range("E4").Value = "QSG3Yn89"

$ws.Range("D5").Value = "yujakov.t.a.2017"
$ws.Range("E5").Value = "X97yaHif"
